# Update "想去人数" (column F) figures across all four sheets to match the
# freshly scraped counts from the gh-pages data refresh (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(5, 6).Value = 5354
$ws.Cells.Item(6, 6).Value = 5355
$ws.Cells.Item(7, 6).Value = 180
$ws.Cells.Item(11, 6).Value = 1197
$ws.Cells.Item(12, 6).Value = 6299
$ws.Cells.Item(14, 6).Value = 78
$ws.Cells.Item(16, 6).Value = 3046
$ws.Cells.Item(17, 6).Value = 265
$ws.Cells.Item(19, 6).Value = 254
$ws.Cells.Item(20, 6).Value = 4008
$ws.Cells.Item(24, 6).Value = 3926
$ws.Cells.Item(25, 6).Value = 189
$ws.Cells.Item(26, 6).Value = 188
$ws.Cells.Item(28, 6).Value = 244
$ws.Cells.Item(29, 6).Value = 255
$ws.Cells.Item(30, 6).Value = 217
$ws.Cells.Item(31, 6).Value = 116
$ws.Cells.Item(36, 6).Value = 28
$ws.Cells.Item(37, 6).Value = 6978
$ws.Cells.Item(38, 6).Value = 31
$ws.Cells.Item(39, 6).Value = 1144
$ws.Cells.Item(40, 6).Value = 552
$ws.Cells.Item(42, 6).Value = 63
$ws.Cells.Item(43, 6).Value = 1418
$ws.Cells.Item(44, 6).Value = 185
$ws.Cells.Item(45, 6).Value = 744
$ws.Cells.Item(46, 6).Value = 2381
$ws.Cells.Item(47, 6).Value = 321
$ws.Cells.Item(49, 6).Value = 793

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(9, 6).Value = 92
$ws.Cells.Item(22, 6).Value = 56
$ws.Cells.Item(25, 6).Value = 830

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 220

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 220
$ws.Cells.Item(8, 6).Value = 5355
$ws.Cells.Item(9, 6).Value = 5355
$ws.Cells.Item(10, 6).Value = 180
$ws.Cells.Item(12, 6).Value = 92
$ws.Cells.Item(15, 6).Value = 1198
$ws.Cells.Item(16, 6).Value = 6299
$ws.Cells.Item(18, 6).Value = 78
$ws.Cells.Item(20, 6).Value = 3047
$ws.Cells.Item(21, 6).Value = 265
$ws.Cells.Item(23, 6).Value = 254
$ws.Cells.Item(24, 6).Value = 4008
$ws.Cells.Item(25, 6).Value = 3926
$ws.Cells.Item(26, 6).Value = 189
$ws.Cells.Item(27, 6).Value = 188
$ws.Cells.Item(28, 6).Value = 244
$ws.Cells.Item(29, 6).Value = 255
$ws.Cells.Item(30, 6).Value = 217
$ws.Cells.Item(31, 6).Value = 116
$ws.Cells.Item(36, 6).Value = 6978
$ws.Cells.Item(37, 6).Value = 31
$ws.Cells.Item(38, 6).Value = 1144
$ws.Cells.Item(39, 6).Value = 552
$ws.Cells.Item(42, 6).Value = 63
$ws.Cells.Item(43, 6).Value = 1418
$ws.Cells.Item(44, 6).Value = 185
$ws.Cells.Item(45, 6).Value = 744
$ws.Cells.Item(46, 6).Value = 2381
$ws.Cells.Item(47, 6).Value = 321
$ws.Cells.Item(48, 6).Value = 793
